# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-123) from 2023-09-02 (serial 45171) to 2023-09-03 (serial 45172).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C123").Value = 45172
